$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.770.58"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("E2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.093.59"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E3").NumberFormat = "General"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E4").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.67"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("E5").NumberFormat = "General"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E6").NumberFormat = "General"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.17"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E7").NumberFormat = "General"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E8").NumberFormat = "General"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("E9").NumberFormat = "General"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E10").NumberFormat = "General"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.066.11"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +27.72%  "
$ws.Range("E11").NumberFormat = "General"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E12").NumberFormat = "General"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.37"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.09%  "
$ws.Range("E13").NumberFormat = "General"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.06"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("E14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.807"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.08%  "
$ws.Range("E15").NumberFormat = "General"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("E16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.096.06"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("E17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.766.21"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("E18").NumberFormat = "General"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.97"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("E19").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.08"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("E20").NumberFormat = "General"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("E21").NumberFormat = "General"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.81"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("E22").NumberFormat = "General"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.37"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("E24").NumberFormat = "General"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("E25").NumberFormat = "General"

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.55"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("E26").NumberFormat = "General"

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.29"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("E27").NumberFormat = "General"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.39%  "
$ws.Range("E28").NumberFormat = "General"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.43"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.95%  "
$ws.Range("E29").NumberFormat = "General"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.29"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("E30").NumberFormat = "General"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.48"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.50%  "
$ws.Range("E31").NumberFormat = "General"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E32").NumberFormat = "General"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.51"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E33").NumberFormat = "General"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.72"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E34").NumberFormat = "General"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("E36").NumberFormat = "General"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("E37").NumberFormat = "General"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("E38").NumberFormat = "General"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E39").NumberFormat = "General"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.00"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E40").NumberFormat = "General"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.90%  "
$ws.Range("E41").NumberFormat = "General"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.01"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("E42").NumberFormat = "General"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.531.98"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("E43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("E44").NumberFormat = "General"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0914"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("E45").NumberFormat = "General"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.13"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("E46").NumberFormat = "General"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.95%  "
$ws.Range("E47").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.13"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("E48").NumberFormat = "General"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("E49").NumberFormat = "General"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("E50").NumberFormat = "General"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.291.21"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.18%  "
$ws.Range("E51").NumberFormat = "General"
